$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("B4").Value = 999536
$ws.Range("C4").Value = 12376
$ws.Range("D4").Value = 137272
$ws.Range("E4").Value = 806084
$ws.Range("G4").Value = 767
$ws.Range("H4").Value = 56180
$ws.Range("B7").Value = 165842
$ws.Range("C7").Value = 3742
$ws.Range("D7").Value = 45513
$ws.Range("E7").Value = 97036
$ws.Range("F7").Value = 4608
$ws.Range("G7").Value = 437
$ws.Range("H7").Value = 23293
$ws.Range("B8").Value = 158213
$ws.Range("C8").Value = 443
$ws.Range("E8").Value = 37692
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = 6021
$ws.Range("A19").Value = "Peru"
$ws.Range("B19").Value = 28699
$ws.Range("C19").Value = 1182
$ws.Range("D19").Value = 8425
$ws.Range("E19").Value = 19492
$ws.Range("F19").Value = 598
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 782
$ws.Range("A20").Value = "India"
$ws.Range("B20").Value = 28380
$ws.Range("C20").Value = 490
$ws.Range("D20").Value = 6523
$ws.Range("E20").Value = 20971
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 886
$ws.Range("E34").Value = 7557
$ws.Range("G34").Value = 22
$ws.Range("H34").Value = 641
$ws.Range("D56").Value = 1140
$ws.Range("E56").Value = 2560
$ws.Range("A86").Value = "Guinea"
$ws.Range("B86").Value = 1163
$ws.Range("C86").Value = 167
$ws.Range("D86").Value = 246
$ws.Range("E86").Value = 910
$ws.Range("H86").Value = 7
$ws.Range("A87").Value = "Costa de Marfil"
$ws.Range("B87").Value = 1150
$ws.Range("D87").Value = 468
$ws.Range("E87").Value = 668
$ws.Range("F87").Value = 0
$ws.Range("H87").Value = 14
$ws.Range("A88").Value = "Hong Kong"
$ws.Range("B88").Value = 1038
$ws.Range("D88").Value = 787
$ws.Range("E88").Value = 247
$ws.Range("F88").Value = 4
$ws.Range("H88").Value = 4
$ws.Range("A89").Value = "Republica de Yibuti"
$ws.Range("B89").Value = 1035
$ws.Range("C89").Value = 12
$ws.Range("D89").Value = 477
$ws.Range("E89").Value = 556
$ws.Range("H89").Value = 2
$ws.Range("B150").Value = 95
$ws.Range("C150").Value = 1
$ws.Range("E150").Value = 49
$ws.Range("A1").Value = "Datos actualizados a 27 de Abril de 2020 a las 20:52"
